$d = $word.ActiveDocument

# Turn on track-changes so the inserted/edited text is recorded as discrete
# runs (rather than being silently coalesced into the existing "test" run).
$d.TrackRevisions = $true

# Append " - test" right after the word "test" (before the trailing
# paragraph mark / _GoBack bookmark).
$tail = $d.Range(4, 4)
$tail.InsertAfter(" - test")

# Turn the leading "t" into "T", which - together with the insertion above -
# leaves the paragraph text reading "Test - test" split across three runs:
# "T", "est", " - test".
$head = $d.Range(0, 1)
$head.Text = "T"

# Accept the tracked changes so the final document has plain runs (no
# <w:ins> markup) while keeping the run boundaries that were introduced.
$d.TrackRevisions = $false
$d.AcceptAllRevisions()
